# feat: new coding for rhythm and music not synch
# also reduce animation time for slot machine

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- "TIPO ATTIVITA'" table (rows 3-7) gains two new columns (C, D) for
#     RITMO ATTIVO / MUSICA ATTIVA, and the existing single numeric code
#     in column B becomes a 3-bit text code split across B/C/D.

# Copy the existing column-B formatting across into the new C/D columns
# for each row of the table before writing values into them.
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B5").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B6").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B7").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row: new sound-state columns (order matches how the strings were
# introduced: D, then C, then B)
$ws.Range("D3").Value = "MUSICA ATTIVA"
$ws.Range("C3").Value = "RITMO ATTIVO"
$ws.Range("B3").Value = "NESSUN SUONO"

# Row 4: base/default state code, now a 3-bit text code
$ws.Range("B4").Value = "_000"
$ws.Range("C4").Value = "_010"
$ws.Range("D4").Value = "_001"

# Row 5: NON COSTANTE
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 110
$ws.Range("D5").Value = 101

# Row 6: OTTO
$ws.Range("B6").Value = 200
$ws.Range("C6").Value = 210
$ws.Range("D6").Value = 201

# Row 7: MOV. ARMONICO
$ws.Range("B7").Value = 300
$ws.Range("C7").Value = 310
$ws.Range("D7").Value = 301

# New column widths for D and E
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668
$ws.Columns.Item(5).ColumnWidth = 12.75

# Update the sheet view: move the selection to C4 (also resets the
# scrolled-down A8 top-left cell from the previous view)
$ws.Activate()
$ws.Range("C4").Select()
